# Weekly CompStat update: new crime data collected (cs-en-us-067pct, week of 5/8/2023-5/14/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Volume/Number and the reporting week date range ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Cells that flip from a numeric count to the text placeholder "0" or "***.*" ---
# (NumberFormat forced to Text so the literal stays a string, then the cells
#  presentation style is restored via a formats-only paste from a donor cell that
#  already carries the correct placeholder style.)
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# --- Cells that flip from the text placeholder back to a real numeric value ---
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100

$excel.CutCopyMode = $false

# --- Plain numeric value updates (new weeks crime counts / recomputed percentages) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -59.375
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = -15.384615384615
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 132
$ws.Range("K16").Value = -26.515151515151
$ws.Range("L16").Value = 36.619718309859
$ws.Range("M16").Value = -38.993710691823
$ws.Range("N16").Value = -89.113355780022
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 243
$ws.Range("J17").Value = 213
$ws.Range("K17").Value = 14.084507042253
$ws.Range("L17").Value = 19.117647058823
$ws.Range("M17").Value = 66.438356164383
$ws.Range("N17").Value = -38.636363636363
$ws.Range("C18").Value = 9
$ws.Range("E18").Value = 80
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 4.761904761904
$ws.Range("I18").Value = 77
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = -6.097560975609
$ws.Range("L18").Value = 75
$ws.Range("M18").Value = -49.673202614379
$ws.Range("N18").Value = -89.622641509434
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -7.8125
$ws.Range("I19").Value = 258
$ws.Range("J19").Value = 266
$ws.Range("K19").Value = -3.007518796992
$ws.Range("L19").Value = 76.712328767123
$ws.Range("M19").Value = 31.632653061224
$ws.Range("N19").Value = -4.089219330855
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 6.25
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = -3.960396039603
$ws.Range("L20").Value = 16.867469879518
$ws.Range("M20").Value = -9.345794392523
$ws.Range("N20").Value = -87.690355329949
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -23.684210526315
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 177
$ws.Range("H21").Value = -1.129943502824
$ws.Range("I21").Value = 792
$ws.Range("J21").Value = 811
$ws.Range("K21").Value = -2.342786683107
$ws.Range("L21").Value = 40.425531914893
$ws.Range("M21").Value = 1.930501930501
$ws.Range("N21").Value = -74.696485623003
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 83.333333333333
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 43.478260869565
$ws.Range("I24").Value = 474
$ws.Range("J24").Value = 331
$ws.Range("K24").Value = 43.202416918429
$ws.Range("L24").Value = 38.59649122807
$ws.Range("M24").Value = 35.042735042735
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -5
$ws.Range("I25").Value = 282
$ws.Range("J25").Value = 242
$ws.Range("K25").Value = 16.528925619834
$ws.Range("L25").Value = 27.601809954751
$ws.Range("M25").Value = -6.622516556291
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -60
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = -6.25
$ws.Range("L26").Value = -28.571428571428
$ws.Range("C27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 30
$ws.Range("L27").Value = 30.434782608695
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = 25
$ws.Range("L28").Value = 17.647058823529
$ws.Range("M28").Value = 17.647058823529
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = 33.333333333333
$ws.Range("L29").Value = -5.882352941176
$ws.Range("M29").Value = 6.666666666666
$ws.Range("N29").Value = -71.929824561403
